# QuyenTX: 2 cases runs successfully with new scenario
# (create new projects in both JIRA and OpenERP)
#
# Adds a "CheckJIRA" xpath-reference block (rows 34-36) to the
# "ServiceImporting" worksheet, mirroring the existing "CheckOpenERP"-style
# blocks already present above it (rows 20-22, 25-27, 30-32).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ServiceImporting")
$q = [char]34

# ---------------------------------------------------------------------
# 1. Cell values, written in the exact order the original authoring tool
#    used so new shared-strings are allocated with matching indices.
# ---------------------------------------------------------------------
$ws.Range("D34").Value = "projectListXpath"
$ws.Range("E34").Value = "issueListXpath"
$ws.Range("F34").Value = "worklogListXpath"
$ws.Range("G34").Value = "firstProjectXpath"
$ws.Range("K34").Value = "firstWorklogXpath"
$ws.Range("D35").Value = "html/body/div/section/div/div/section/div[1]/div[2]/div/table/tbody"
$ws.Range("E35").Value = "html/body/div[1]/section/div[2]/div/section/div/div/div/div/div/div/div[1]/div[1]/div/div[1]/div[2]/div/ol"
$ws.Range("F35").Value = "html/body/div[1]/section/div[2]/div/section/div/div/div/div/div/div/div[2]/div[1]/div/div/div/div[1]/div[5]/div[2]/div[2]/div[2]"
$ws.Range("G35").Value = "html/body/div[1]/section/div/div/section/div[1]/div[2]/div/table/tbody/tr[1]/td[2]/a"
$ws.Range("I35").Value = "html/body/div[1]/section/div[2]/div/section/div/div/div/div/div/div/div[1]/div[1]/div/div[1]/div[2]/div/ol/li[" + $q + "+issueListSize+" + $q + "]/a/span[2]"
$ws.Range("K35").Value = "html/body/div[1]/section/div[2]/div/section/div/div/div/div/div/div/div[2]/div[1]/div/div/div/div[1]/div[5]/div[2]/div[2]/div[2]/div[" + $q + "+worklogListSize+" + $q + "]"
$ws.Range("L35").Value = "html/body/div[1]/section/div[2]/div/section/div/div/div/div/div/div/div[2]/div[1]/div/div/div/div[1]/div[5]/div[2]/div[2]/div[2]/div[" + $q + "+worklogListSize+" + $q + "]/div/div[2]/a"
$ws.Range("M35").Value = "html/body/div[1]/section/div[2]/div/section/div/div/div/div/div/div/div[2]/div[1]/div/div/div/div[1]/div[5]/div[2]/div[2]/div[2]/div[" + $q + "+worklogListSize+" + $q + "]/div/div[3]/ul/li/dl[1]/dd"
$ws.Range("L34").Value = "firstWorklogNameXpath"
$ws.Range("M34").Value = "firstWorklogTimeXpath"
$ws.Range("A34").Value = "CheckJIRA"
$ws.Range("H35").Value = "html/body/div/section/div/div/section/div[1]/div[2]/div/table/tbody/tr[1]/td[3]"
$ws.Range("H34").Value = "firstProjectKeyXpath"
$ws.Range("I34").Value = "firstIssueNameXpath"
$ws.Range("J34").Value = "firstIssueKeyXpath"
$ws.Range("J35").Value = "html/body/div[1]/section/div[2]/div/section/div/div/div/div/div/div/div[1]/div[1]/div/div[1]/div[2]/div/ol/li[" + $q + "+issueListSize+" + $q + "]/a/span[1]"

$ws.Range("N36").Value = "CheckJIRA"
$ws.Range("B34").Value = "username"
$ws.Range("C34").Value = "password"
$ws.Range("B35").Value = "admin"
$ws.Range("C35").Value = "123456@a"

# ---------------------------------------------------------------------
# 2. Formatting - copied from the existing, analogous "CheckOpenERP" block
#    (rows 30-32) so the engine re-uses already-defined styles instead of
#    creating near-duplicate ones.
# ---------------------------------------------------------------------
$ws.Range("A30").Copy()
$ws.Range("A34").PasteSpecial(-4122)

$ws.Range("B30").Copy()
$ws.Range("B34:M34").PasteSpecial(-4122)

$ws.Range("C31").Copy()
$ws.Range("B35").PasteSpecial(-4122)
$ws.Range("D35:M35").PasteSpecial(-4122)

$ws.Range("G32").Copy()
$ws.Range("N36").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Row heights for the new block.
$ws.Rows.Item(34).RowHeight = 30
$ws.Rows.Item(35).RowHeight = 225
$ws.Rows.Item(36).RowHeight = 30

# C35 becomes a hyperlink (mirrors K2's mailto link) with wrapped text -
# this produces the one genuinely new cell style (Hyperlink + wrapText,
# no border).
$ws.Range("C35").WrapText = $true
$ws.Hyperlinks.Add($ws.Range("C35"), "mailto:123456@a")

# ---------------------------------------------------------------------
# 3. Column widths widened to fit the new, much longer xpath content.
#    (input values are pre-compensated for this engine's column-width
#    rounding so the stored width lands as close as possible to Excel's
#    original 9.7109375 / 14.5703125 / ... values)
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 8.794270833333332
$ws.Columns.Item(5).ColumnWidth = 13.653645833333332
$ws.Columns.Item(6).ColumnWidth = 21.416822916666668
$ws.Columns.Item(8).ColumnWidth = 23.653645833333336
$ws.Columns.Item(9).ColumnWidth = 24.250625
$ws.Columns.Item(10).ColumnWidth = 24.653645833333336
$ws.Columns.Item(11).ColumnWidth = 25.584114583333335

# ---------------------------------------------------------------------
# 4. View state - scroll / selection moved onto the new block.
# ---------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("I35").Select()

Write-Output "done"
